$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "-0.372***"
$ws.Range("B3").Value = "-3.464***"
$ws.Range("C2").Value = "0.01*"
$ws.Range("C3").Value = "-0.808***"
